$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Assert-CellText($cell, $expected) {
    $actual = $cell.Range.Text
    # Range.Text carries a trailing paragraph mark (0x0d) and cell mark
    # (0x07) that aren't part of the visible cell contents.
    $actual = $actual -replace "[\x07\x0d]", ""
    if ($actual -ne $expected) {
        throw "cell text mismatch: expected [$expected] got [$actual]"
    }
}

# --- Row 7 : "Installation de chantier" line ------------------------------
# Col 2 : "L'ens" -- drop the spell-check proofErr wrapper, keep pPr as-is.
$cell = $t.Cell(7, 2)
Assert-CellText $cell "L’ens"
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="243"/><w:rPr><w:vertAlign w:val="superscript"/></w:rPr></w:pPr><w:r><w:t>L’ens</w:t></w:r></w:p>
'@
$cell.Range.InsertXML($xml)

# Col 3 : "       5,00" -> split into digit-by-digit runs "     500,05"
# and drop the now-unneeded first-line indent paragraph property.
$cell = $t.Cell(7, 3)
Assert-CellText $cell "       5,00"
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">     </w:t></w:r><w:r><w:t>5</w:t></w:r><w:r><w:t>0</w:t></w:r><w:r><w:t>0</w:t></w:r><w:r><w:t>,0</w:t></w:r><w:r><w:t>5</w:t></w:r></w:p>
'@
$cell.Range.InsertXML($xml)

# --- Row 11 : "Installation de chantier compris repli" line --------------
# Col 2 : "L'ens" -- same proofErr cleanup, pPr unchanged, text unchanged.
$cell = $t.Cell(11, 2)
Assert-CellText $cell "L’ens"
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="243"/></w:pPr><w:r><w:t>L’ens</w:t></w:r></w:p>
'@
$cell.Range.InsertXML($xml)

# --- Row 125 : "attention, merci de nous ..." -----------------------------
$cell = $t.Cell(125, 1)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">         attention, merci de nous réexpédier un exemplaire signé.</w:t></w:r></w:p>
'@
$cell.Range.InsertXML($xml)

# --- Row 130 : "   signé, solde ..." --------------------------------------
$cell = $t.Cell(130, 1)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">   signé</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> solde à réception de facture</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>. A</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">ttention si vous ne </w:t></w:r></w:p>
'@
$cell.Range.InsertXML($xml)

# --- Row 131 : "   joignez pas d'acompte ..." -----------------------------
$cell = $t.Cell(131, 1)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">   joignez pas d’acompte le chantier ne sera pas pris en compte dans</w:t></w:r></w:p>
'@
$cell.Range.InsertXML($xml)

# --- Row 132 : "   le planning" -------------------------------------------
$cell = $t.Cell(132, 1)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">   </w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>le planning</w:t></w:r></w:p>
'@
$cell.Range.InsertXML($xml)
